$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (row 19) since the roster now has 17 players instead of 18
$ws.Rows(19).Delete()

# Update player roster data in A2:C18 with the new roster
$ws.Cells.Item(2,1).Value = "Chris Paul"
$ws.Cells.Item(2,2).Value = "PG"
$ws.Cells.Item(2,3).Value = "San Antonio Spurs"
$ws.Cells.Item(3,1).Value = "Jalen Green"
$ws.Cells.Item(3,2).Value = "PG,SG"
$ws.Cells.Item(3,3).Value = "Houston Rockets"
$ws.Cells.Item(4,1).Value = "Buddy Hield"
$ws.Cells.Item(4,2).Value = "SG,SF"
$ws.Cells.Item(4,3).Value = "Golden State Warriors"
$ws.Cells.Item(5,1).Value = "Rudy Gobert"
$ws.Cells.Item(5,2).Value = "C"
$ws.Cells.Item(5,3).Value = "Minnesota Timberwolves"
$ws.Cells.Item(6,1).Value = "Jakob Poeltl"
$ws.Cells.Item(6,2).Value = "C"
$ws.Cells.Item(6,3).Value = "Toronto Raptors"
$ws.Cells.Item(7,1).Value = "Trayce Jackson-Davis"
$ws.Cells.Item(7,2).Value = "PF,C"
$ws.Cells.Item(7,3).Value = "Golden State Warriors"
$ws.Cells.Item(8,1).Value = "Khris Middleton"
$ws.Cells.Item(8,2).Value = "SF"
$ws.Cells.Item(8,3).Value = "Milwaukee Bucks"
$ws.Cells.Item(9,1).Value = "Bennedict Mathurin"
$ws.Cells.Item(9,2).Value = "SG,SF"
$ws.Cells.Item(9,3).Value = "Indiana Pacers"
$ws.Cells.Item(10,1).Value = "Pascal Siakam"
$ws.Cells.Item(10,2).Value = "SF,PF,C"
$ws.Cells.Item(10,3).Value = "Indiana Pacers"
$ws.Cells.Item(11,1).Value = "Paolo Banchero"
$ws.Cells.Item(11,2).Value = "SF,PF"
$ws.Cells.Item(11,3).Value = "Orlando Magic"
$ws.Cells.Item(12,1).Value = "Jaylen Brown"
$ws.Cells.Item(12,2).Value = "SG,SF"
$ws.Cells.Item(12,3).Value = "Boston Celtics"
$ws.Cells.Item(13,1).Value = "Nikola Jokic"
$ws.Cells.Item(13,2).Value = "C"
$ws.Cells.Item(13,3).Value = "Denver Nuggets"
$ws.Cells.Item(14,1).Value = "Russell Westbrook"
$ws.Cells.Item(14,2).Value = "PG,SG"
$ws.Cells.Item(14,3).Value = "Denver Nuggets"
$ws.Cells.Item(15,1).Value = "Dejounte Murray"
$ws.Cells.Item(15,2).Value = "PG,SG"
$ws.Cells.Item(15,3).Value = "New Orleans Pelicans"
$ws.Cells.Item(16,1).Value = "Deni Avdija"
$ws.Cells.Item(16,2).Value = "SF,PF"
$ws.Cells.Item(16,3).Value = "Portland Trail Blazers"
$ws.Cells.Item(17,1).Value = "Chet Holmgren"
$ws.Cells.Item(17,2).Value = "PF,C"
$ws.Cells.Item(17,3).Value = "Oklahoma City Thunder"
$ws.Cells.Item(18,1).Value = "Jalen Suggs"
$ws.Cells.Item(18,2).Value = "PG,SG"
$ws.Cells.Item(18,3).Value = "Orlando Magic"
